$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(4)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Paragraph 2: "One of my proudest tech moments was to make CWeave and CWeb
# (and LaTeX) run in a Windows computer circa 1998."
# Split the trailing run so "LaTeX" becomes its own run in Consolas.
$para2 = $tr.Paragraphs(2)
$lastRun = $para2.Runs(5)
$lastRun.Text = " (and "
$null = $lastRun.InsertAfter("LaTeX")
$null = $para2.Runs(6).InsertAfter(") run in a Windows computer circa 1998.")
$para2.Runs(6).Font.Name = "Consolas"

# Paragraph 3: replace the Jupyter/Quarto sentence with the new one, keeping
# it as a single run.
$para3 = $tr.Paragraphs(3)
$para3.Runs(1).Text = "Jupyter implements the literate programming paradigm, but I haven’t seen the Markdown part gain a lot of traction."
